$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "['MEC-3B-Cont.Lóg.Prog CLP', -, -, -]"
$ws.Range("E2").Value = "-"
$ws.Range("F2").Value = "-"

$ws.Range("E3").Value = "['MEC-2B-Des. Maq. Cad._T2', -]"
$ws.Range("F3").Value = "-"

$ws.Range("C4").Value = "[-, 'MEC-2B-Des. Maq. Cad._T1']"
$ws.Range("E4").Value = "['MEC-2B-Des. Maq. Cad._T2', -]"
$ws.Range("F4").Value = "-"

$ws.Range("C6").Value = "[-, 'MEC-2B-Des. Maq. Cad._T1']"
$ws.Range("D6").Value = "['MEC-3B-Cont.Lóg.Prog CLP', -, -, -]"
$ws.Range("E6").Value = "['MEC-2B-Des. Maq. Cad._T2', -]"

$ws.Range("C7").Value = "[-, 'MEC-2B-Des. Maq. Cad._T1']"
$ws.Range("D7").Value = "['MEC-3B-Cont.Lóg.Prog CLP', -, -, -]"
$ws.Range("F7").Value = "-"

$ws.Range("D8").Value = "['MEC-3B-Cont.Lóg.Prog CLP', -, -, -]"
$ws.Range("E8").Value = "-"
$ws.Range("F8").Value = "-"

$ws.Range("B12").Value = "[-, 'MEC-1A-Desenho tecnico mecanico']"

$ws.Range("B14").Value = "[-, 'MEC-1A-Desenho tecnico mecanico']"
$ws.Range("F14").Value = "-"

$ws.Range("B15").Value = "[-, 'MEC-1A-Desenho tecnico mecanico']"
$ws.Range("F15").Value = "-"

$ws.Range("F16").Value = "-"

$ws.Range("B18").Value = "-"
$ws.Range("C18").Value = "-"
$ws.Range("D18").Value = "[-, 'MEC-1NA-Desenho tecnico mecanico – T2']"

$ws.Range("B19").Value = "-"
$ws.Range("C19").Value = "-"

$ws.Range("B20").Value = "[-, -, -, 'MEC-2NB-C.pneumática']"
$ws.Range("D20").Value = "[-, 'MEC-1NA-Desenho tecnico mecanico – T2']"

$ws.Range("B21").Value = "['MEC-2NB-C.pneumática', -, -, -]"
$ws.Range("D21").Value = "[-, 'MEC-1NA-Desenho tecnico mecanico – T2']"
